$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buffer")
$ws.Activate()
Write-Output ("before delete AutoFilterMode=" + $ws.AutoFilterMode)
$ws.Columns.Item(1).Delete()
Write-Output ("after delete AutoFilterMode=" + $ws.AutoFilterMode)
Write-Output ("after delete ref=" + $ws.AutoFilter.Range.Address())
$ws.AutoFilterMode = $false
Write-Output ("after clear mode=" + $ws.AutoFilterMode)
$ws.Range("A1:C1048576").AutoFilter()
Write-Output ("after AutoFilter() mode=" + $ws.AutoFilterMode)
try {
  Write-Output ("after AutoFilter() ref=" + $ws.AutoFilter.Range.Address())
} catch {
  Write-Output ("error getting ref: " + $_)
}
